$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Reference cell with the default (unstyled) cell style, used to restore
# the style after a temporary text-number-format trick below.
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "26.959.82"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.554.09"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "206.86"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "22.03"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "0.0595"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.773.11"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "1.547.24"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "26.923.04"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "61.65"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "217.45"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "153.39"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "6.64"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0470"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "3.12"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("D34").Value = "1.418.66"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "1.09"
$ws.Range("E35").Value = "  +13.32%  "
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D39").Value = "0.526"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "64.48"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "1.687.55"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "87.24"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "0.0₇0998"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  +0.51%  "
